# Auto-generated edit script
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

$ws.Range("B3").Value = 'operationHttpStatusMapper'
$ws.Range("D3").Value = 'org.andante.mappers.OperationHttpStatusMapper'
$ws.Range("B4").Value = 'IDENTIFIERS_LIST_NULL_MESSAGE'
$ws.Range("D4").Value = 'java.lang.String'
$ws.Range("B5").Value = 'NEGATIVE_PAGE_ERROR_MESSAGE'
$ws.Range("B6").Value = 'ACTIVITY_EMAIL_NOT_VALID_MESSAGE'
$ws.Range("B8").Value = 'activityDTOModelMapper'
$ws.Range("D8").Value = 'org.andante.activity.controller.mapper.ActivityDTOModelMapper'
$ws.Range("B9").Value = 'NULL_PAGE_SIZE_ERROR_MESSAGE'
$ws.Range("B10").Value = 'NON_POSITIVE_PAGE_SIZE_MESSAGE'
$ws.Range("D10").Value = 'java.lang.String'
$ws.Range("B11").Value = 'ACTIVITY_IDENTIFIER_NOT_BLANK_MESSAGE'
$ws.Range("B12").Value = 'IDENTIFIERS_LIST_MESSAGE'
$ws.Range("B13").Value = 'ACTIVITY_EMAIL_BLANK_MESSAGE'
$ws.Range("B18").Value = 'affectedUsers'
$ws.Range("D18").Value = 'java.util.Set'
$ws.Range("B19").Value = 'acknowledgedUsers'
$ws.Range("D19").Value = 'java.util.Set'
$ws.Range("B20").Value = 'description'
$ws.Range("D20").Value = 'java.lang.String'
$ws.Range("B21").Value = 'domain'
$ws.Range("D21").Value = 'org.andante.activity.enums.Domain'
$ws.Range("B22").Value = 'relatedId'
$ws.Range("D22").Value = 'java.lang.String'
$ws.Range("B23").Value = 'eventTimestamp'
$ws.Range("D23").Value = 'java.time.LocalDateTime'
$ws.Range("B24").Value = 'id'
$ws.Range("D24").Value = 'java.lang.String'
$ws.Range("B25").Value = 'priority'
$ws.Range("D25").Value = 'org.andante.activity.enums.Priority'
$ws.Range("B26").Value = 'affectsAll'
$ws.Range("D26").Value = 'java.lang.Boolean'
$ws.Range("B31").Value = 'key'
$ws.Range("B32").Value = 'observedUsers'
$ws.Range("D32").Value = 'java.util.Set'
$ws.Range("B34").Value = 'username'
$ws.Range("D34").Value = 'java.lang.String'
$ws.Range("B35").Value = 'imageUrl'
$ws.Range("B44").Value = 'username'
$ws.Range("B47").Value = 'password'
$ws.Range("B48").Value = 'host'
$ws.Range("B49").Value = 'USERNAME_NULL_ERROR_MESSAGE'
$ws.Range("D49").Value = 'java.lang.String'
$ws.Range("B50").Value = 'IDENTIFIERS_LIST_NULL_ERROR_MESSAGE'
$ws.Range("B52").Value = 'profileService'
$ws.Range("D52").Value = 'org.andante.activity.logic.ProfileService'
$ws.Range("B53").Value = 'userProfileService'
$ws.Range("D53").Value = 'org.andante.activity.logic.UserProfileService'
$ws.Range("B55").Value = 'IDENTIFIERS_LIST_SIZE_ERROR_MESSAGE'
$ws.Range("D55").Value = 'java.lang.String'
$ws.Range("B56").Value = 'userProfileMapper'
$ws.Range("D56").Value = 'org.andante.activity.controller.mapper.UserProfileDTOModelMapper'
$ws.Range("B57").Value = 'IMAGE_URL_BLANK_ERROR_MESSAGE'
$ws.Range("D57").Value = 'java.lang.String'
$ws.Range("B59").Value = 'key'
$ws.Range("D59").Value = 'java.lang.String'
$ws.Range("B61").Value = 'username'
$ws.Range("B62").Value = 'observingUsers'
$ws.Range("B64").Value = 'observedUsers'
$ws.Range("D64").Value = 'java.util.Set'
$ws.Range("B66").Value = 'domain'
$ws.Range("D66").Value = 'org.andante.activity.enums.Domain'
$ws.Range("B67").Value = 'affectedUsers'
$ws.Range("D67").Value = 'java.util.Set'
$ws.Range("B68").Value = 'affectsAll$value'
$ws.Range("D68").Value = 'java.lang.Boolean'
$ws.Range("B69").Value = 'eventTimestamp'
$ws.Range("D69").Value = 'java.time.LocalDateTime'
$ws.Range("B70").Value = 'affectsAll$set'
$ws.Range("D70").Value = 'boolean'
$ws.Range("B72").Value = 'priority'
$ws.Range("D72").Value = 'org.andante.activity.enums.Priority'
$ws.Range("B73").Value = 'relatedId'
$ws.Range("D73").Value = 'java.lang.String'
$ws.Range("B75").Value = 'acknowledgedUsers'
$ws.Range("D75").Value = 'java.util.Set'
$ws.Range("B82").Value = 'keycloakGetUserPath'
$ws.Range("B83").Value = 'userProfileService'
$ws.Range("D83").Value = 'org.andante.activity.logic.UserProfileService'
$ws.Range("B85").Value = 'keycloakAdminTokenPath'
$ws.Range("B86").Value = 'adminUsername'
$ws.Range("B87").Value = 'adminPassword'
$ws.Range("D87").Value = 'java.lang.String'
$ws.Range("B88").Value = 'emailAddress'
$ws.Range("D88").Value = 'java.lang.String'
$ws.Range("B89").Value = 'subscriptionDate'
$ws.Range("D89").Value = 'java.time.LocalDateTime'
$ws.Range("B90").Value = 'isConfirmed'
$ws.Range("D90").Value = 'java.lang.Boolean'
$ws.Range("B91").Value = 'subscriptionDate'
$ws.Range("D91").Value = 'java.time.LocalDateTime'
$ws.Range("B93").Value = 'emailAddress'
$ws.Range("D93").Value = 'java.lang.String'
$ws.Range("B95").Value = 'observed'
$ws.Range("D95").Value = 'java.util.Set'
$ws.Range("B97").Value = 'imageUrl'
$ws.Range("B98").Value = 'observers'
$ws.Range("D98").Value = 'java.util.Set'
$ws.Range("B99").Value = 'username'
$ws.Range("D99").Value = 'java.lang.String'
$ws.Range("B100").Value = 'id'
$ws.Range("D100").Value = 'java.lang.String'
$ws.Range("B103").Value = 'logo'
$ws.Range("D103").Value = 'org.springframework.core.io.Resource'
$ws.Range("B104").Value = 'NEWSLETTER_TEMPLATE'
$ws.Range("D104").Value = 'java.lang.String'
$ws.Range("B106").Value = 'sender'
$ws.Range("D106").Value = 'java.lang.String'
$ws.Range("B107").Value = 'mailSender'
$ws.Range("D107").Value = 'org.springframework.mail.javamail.JavaMailSender'
$ws.Range("B108").Value = 'templateEngine'
$ws.Range("D108").Value = 'org.thymeleaf.TemplateEngine'
$ws.Range("B110").Value = 'USER_NOT_FOUND_EXCEPTION_MESSAGE'
$ws.Range("B111").Value = 'userProfileRepository'
$ws.Range("D111").Value = 'org.andante.activity.repository.UserProfileRepository'
$ws.Range("B112").Value = 'USER_CONFLICT_EXCEPTION_MESSAGE'
$ws.Range("B113").Value = 'userProfileModelEntityMapper'
$ws.Range("D113").Value = 'org.andante.activity.logic.mapper.UserProfileModelEntityMapper'
$ws.Range("B115").Value = 'NEWSLETTER_NOT_FOUND_EXCEPTION_MESSAGE'
$ws.Range("D115").Value = 'java.lang.String'
$ws.Range("B116").Value = 'newsletterRepository'
$ws.Range("D116").Value = 'org.andante.activity.repository.NewsletterRepository'
$ws.Range("B117").Value = 'newsletterMapper'
$ws.Range("D117").Value = 'org.andante.activity.logic.mapper.NewsletterModelEntityMapper'
$ws.Range("B118").Value = 'NEWSLETTER_CONFLICT_EXCEPTION_MESSAGE'
$ws.Range("D118").Value = 'java.lang.String'
$ws.Range("B119").Value = 'id'
$ws.Range("D119").Value = 'java.lang.String'
$ws.Range("B120").Value = 'imageUrl'
$ws.Range("D120").Value = 'java.lang.String'
$ws.Range("B121").Value = 'observers'
$ws.Range("D121").Value = 'java.util.Set'
$ws.Range("B122").Value = 'username'
$ws.Range("B123").Value = 'observed'
$ws.Range("D123").Value = 'java.util.Set'
$ws.Range("B129").Value = 'privateToken'
$ws.Range("B130").Value = 'databaseId'
$ws.Range("B132").Value = 'USER_NOT_AFFECTED_EXCEPTION_MESSAGE'
$ws.Range("B133").Value = 'activityModelEntityMapper'
$ws.Range("D133").Value = 'org.andante.activity.logic.mapper.ActivityModelEntityMapper'
$ws.Range("B134").Value = 'ACTIVITY_CONFLICT_EXCEPTION_MESSAGE'
$ws.Range("D134").Value = 'java.lang.String'
$ws.Range("B135").Value = 'ACTIVITY_NOT_FOUND_EXCEPTION_MESSAGE'
$ws.Range("B136").Value = 'rsqlParser'
$ws.Range("D136").Value = 'cz.jirutka.rsql.parser.RSQLParser'
$ws.Range("B137").Value = 'rsqlVisitor'
$ws.Range("D137").Value = 'org.andante.rsql.PersistentRSQLVisitor'
$ws.Range("B138").Value = 'activityRepository'
$ws.Range("D138").Value = 'org.andante.activity.repository.ActivityRepository'
$ws.Range("B140").Value = 'relatedId'
$ws.Range("D140").Value = 'java.lang.String'
$ws.Range("B141").Value = 'priority'
$ws.Range("D141").Value = 'org.andante.activity.enums.Priority'
$ws.Range("B142").Value = 'domain'
$ws.Range("D142").Value = 'org.andante.activity.enums.Domain'
$ws.Range("B143").Value = 'description'
$ws.Range("D143").Value = 'java.lang.String'
$ws.Range("B144").Value = 'id'
$ws.Range("D144").Value = 'java.lang.String'
$ws.Range("B145").Value = 'acknowledgedUsers'
$ws.Range("D145").Value = 'java.util.Set'
$ws.Range("B146").Value = 'eventTimestamp'
$ws.Range("D146").Value = 'java.time.LocalDateTime'
$ws.Range("B147").Value = 'affectedUsers'
$ws.Range("D147").Value = 'java.util.Set'
$ws.Range("B148").Value = 'affectsAll'
$ws.Range("D148").Value = 'java.lang.Boolean'
$ws.Range("B150").Value = 'emailSender'
$ws.Range("D150").Value = 'org.andante.activity.controller.email.EmailSender'
$ws.Range("B151").Value = 'newsletterService'
$ws.Range("D151").Value = 'org.andante.activity.logic.NewsletterService'
$ws.Range("B154").Value = 'EMAIL_NOT_VALID_ERROR_MESSAGE'
$ws.Range("D154").Value = 'java.lang.String'
$ws.Range("B158").Value = 'emailAddress'
$ws.Range("D158").Value = 'java.lang.String'
$ws.Range("B159").Value = 'isConfirmed'
$ws.Range("D159").Value = 'java.lang.Boolean'
$ws.Range("B160").Value = 'subscriptionDate'
$ws.Range("D160").Value = 'java.time.LocalDateTime'
$ws.Range("B164").Value = 'affectedUsers'
$ws.Range("D164").Value = 'java.util.Set'
$ws.Range("B165").Value = 'domain'
$ws.Range("D165").Value = 'org.andante.activity.enums.Domain'
$ws.Range("B166").Value = 'acknowledgedUsers'
$ws.Range("D166").Value = 'java.util.Set'
$ws.Range("B167").Value = 'relatedId'
$ws.Range("D167").Value = 'java.lang.String'
$ws.Range("B168").Value = 'description'
$ws.Range("B170").Value = 'key'
$ws.Range("D170").Value = 'java.lang.String'
$ws.Range("B171").Value = 'affectsAll'
$ws.Range("D171").Value = 'java.lang.Boolean'
$ws.Range("B173").Value = 'entrySet'
$ws.Range("D173").Value = 'java.util.Set'
$ws.Range("B174").Value = 'DEFAULT_LOAD_FACTOR'
$ws.Range("D174").Value = 'float'
$ws.Range("B175").Value = 'this$0'
$ws.Range("D175").Value = 'org.andante.activity.logic.impl.DefaultRecommendationService'
$ws.Range("B177").Value = 'keySet'
$ws.Range("D177").Value = 'java.util.Set'
$ws.Range("B178").Value = 'UNTREEIFY_THRESHOLD'
$ws.Range("B179").Value = 'DEFAULT_INITIAL_CAPACITY'
$ws.Range("D179").Value = 'int'
$ws.Range("B180").Value = 'MAXIMUM_CAPACITY'
$ws.Range("D180").Value = 'int'
$ws.Range("B181").Value = 'val$productOutputDTO'
$ws.Range("D181").Value = 'org.andante.product.dto.ProductOutputDTO'
$ws.Range("B182").Value = 'loadFactor'
$ws.Range("D182").Value = 'float'
$ws.Range("B184").Value = 'size'
$ws.Range("D184").Value = 'int'
$ws.Range("B185").Value = 'TREEIFY_THRESHOLD'
$ws.Range("B186").Value = 'MIN_TREEIFY_CAPACITY'
$ws.Range("D186").Value = 'int'
$ws.Range("B187").Value = 'table'
$ws.Range("D187").Value = 'java.util.HashMap$Node[]'
$ws.Range("B188").Value = 'values'
$ws.Range("D188").Value = 'java.util.Collection'
